$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price / Volume updates for rows whose coin identity did not change ---
$updates = @{
    2  = @("28.224.67", "  -0.65%  ")
    3  = @("1.809.34",  "  +0.90%  ")
    4  = @("0.9988",    "  -0.60%  ")
    5  = @("316.83",    "  +0.93%  ")
    6  = @("0.9985",    "  -0.48%  ")
    7  = @("0.5337",    "  +0.95%  ")
    8  = @("0.3789",    "  +0.07%  ")
    9  = @("0.07482",   "  -0.46%  ")
    10 = @("41.98",     "  -1.98%  ")
    11 = @("1.102",     "  -1.06%  ")
    12 = @("0.9983",    "  -0.67%  ")
    13 = @("20.66",     "  -1.27%  ")
    14 = @("6.179",     "  +0.10%  ")
    17 = @("89.57",     "  -0.52%  ")
    18 = @("0.00001064","  +0.05%  ")
    19 = @("0.06500",   "  +0.65%  ")
    20 = @("0.9985",    "  -0.44%  ")
    21 = @("17.29",     "  +0.22%  ")
    22 = @("5.912",     "  -0.26%  ")
    23 = @("28.250.50", "  -0.62%  ")
    24 = @("11.20",     "  -1.22%  ")
    25 = @("2.087",     "  -2.33%  ")
    26 = @("156.25",    "  -2.35%  ")
    27 = @("20.45",     "  -0.35%  ")
    28 = @("2.010.81",  "  +0.65%  ")
    29 = @("2.326",     "  -2.44%  ")
    30 = @("121.89",    "  -0.82%  ")
    31 = @("1.122",     "  +1.17%  ")
    32 = @("0.1080",    "  +7.04%  ")
    37 = @("0.02291",   "  -0.98%  ")
    38 = @("5.067",     "  -0.87%  ")
    39 = @("8.501",     "  -2.17%  ")
    40 = @("0.6196",    "  -1.76%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}

# Row 41 only has the Volume(1h) column changed (Price stays "11.16")
$ws.Cells.Item(41, 5).Value = "  -3.17%  "

# --- Rows whose coin identity (Coin/Link) also changed, along with Price/Volume ---
$rowData = @{
    15 = @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "7.334", "  -0.23%  ")
    16 = @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.802.21", "  +0.63%  ")
    33 = @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "5.608", "  -1.54%  ")
    34 = @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.07390", "  +11.82%  ")
    35 = @("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "3.650", "  -0.73%  ")
    36 = @("Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.2240", "  -2.77%  ")
    42 = @("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.184", "  -2.03%  ")
    43 = @("WEMIXTOKEN", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "1.433", "  +2.61%  ")
    44 = @("Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "0.9975", "  -0.58%  ")
    45 = @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "13.36", "  -1.13%  ")
    46 = @("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "3.674", "  +0.19%  ")
    47 = @("Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.5794", "  -2.38%  ")
    48 = @("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "125.34", "  -0.11%  ")
    49 = @("EOS", "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos", "1.191", "  +2.98%  ")
    50 = @("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "1.935", "  -2.43%  ")
    51 = @("Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.06851", "  -1.12%  ")
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
}
